$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.059.77'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.679.65'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '215.75'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").Value = '  +5.75%  '
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").Value = '1.917.29'
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = '1.640.82'
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").Value = '66.33'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '27.054.65'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").Value = '8.15'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").Value = '236.44'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '4.47'
$ws.Range("E22").Value = '  +2.46%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").Value = '147.25'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").Value = '7.27'
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").Value = '  +3.83%  '
$ws.Range("E28").Value = '  -1.51%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '1.546.38'
$ws.Range("E33").Value = '  +6.11%  '
$ws.Range("E35").Value = '  +5.20%  '
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("D37").Value = '0.588'
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").Value = '0.915'
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("E40").Value = '  +6.93%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '67.84'
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("D44").Value = '2.26'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").Value = '1.822.22'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").Value = '90.51'
$ws.Range("E48").Value = '  +2.66%  '
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").Value = '8.05'
$ws.Range("E51").Value = '  +7.23%  '
